# Apply updates to "Hortaliza, Vega Modelo de Temuco - Ramas de apio" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: row number -> hashtable of column letter -> new value
$updates = @{
    2  = @{ D = 44680; J = 20 }
    3  = @{ D = 44315 }
    4  = @{ D = 44749; J = 65; K = 6000; L = 6000; M = 6000; P = 6000 }
    5  = @{ D = 44777; J = 25; K = 5000; L = 5000; M = 5000; P = 5000 }
    6  = @{ D = 44291; J = 35 }
    7  = @{ D = 44509; J = 20; K = 4000; L = 4000; M = 4000; P = 4000 }
    8  = @{ D = 44280; J = 55 }
    9  = @{ D = 44649; J = 20 }
    10 = @{ D = 44957; J = 20; K = 5000; L = 5000; M = 5000; P = 5000 }
    11 = @{ D = 44956; J = 40 }
    12 = @{ D = 44301; J = 40; K = 3000; L = 3000; M = 3000; P = 3000 }
    13 = @{ D = 44781; J = 40; K = 5000; L = 5000; M = 5000; P = 5000 }
    14 = @{ D = 44312; J = 50; K = 4000; L = 4000; M = 4000; P = 4000 }
    16 = @{ D = 44504; J = 55 }
    17 = @{ D = 44259; J = 30 }
    18 = @{ D = 44313; K = 4000; L = 4000; M = 4000; P = 4000 }
    19 = @{ D = 44656; J = 85; K = 5000; L = 5000; M = 5000; P = 5000 }
    20 = @{ D = 44497; J = 20; K = 4000; L = 4000; M = 4000; P = 4000 }
    21 = @{ D = 44365; J = 55; K = 5000; L = 5000; M = 5000; P = 5000 }
    22 = @{ D = 44176; J = 10 }
    23 = @{ D = 44498; J = 40 }
    24 = @{ D = 44966; J = 40; K = 5000; L = 5000; M = 5000; P = 5000 }
    25 = @{ D = 44316; K = 4000; L = 4000; M = 4000; P = 4000 }
    26 = @{ D = 44679; J = 50 }
    27 = @{ D = 44959; J = 40 }
    28 = @{ D = 44508; J = 30; K = 4000; L = 4000; M = 4000; P = 4000 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $cols[$col]
    }
}
